$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8815722556963124
$ws.Range("C2").Value = 0.2114409464772109
$ws.Range("D2").Value = 0.5999240253903508
$ws.Range("E2").Value = 0.2300334483836153
$ws.Range("G2").Value = 0.5707739732063359
$ws.Range("H2").Value = 0.7025952422588446
$ws.Range("I2").Value = 0.4993916693081673
$ws.Range("J2").Value = 0.1106633300735709
$ws.Range("M2").Value = 0.4233633840258975
$ws.Range("O2").Value = 2.523295925708524
$ws.Range("B3").Value = 0.779178266181475
$ws.Range("C3").Value = 0.1848042638789309
$ws.Range("D3").Value = 0.5968583030539065
$ws.Range("E3").Value = 0.2303273256448577
$ws.Range("G3").Value = 0.5744778696322186
$ws.Range("H3").Value = 0.7090904665238043
$ws.Range("I3").Value = 0.5086765228882459
$ws.Range("J3").Value = 0.1116576869621291
$ws.Range("M3").Value = 0.3941726583459868
$ws.Range("O3").Value = 2.544431225188191
$ws.Range("B4").Value = 0.7161530286056177
$ws.Range("C4").Value = 0.1683884175627384
$ws.Range("D4").Value = 0.5952790985280387
$ws.Range("E4").Value = 0.2306272622248571
$ws.Range("G4").Value = 0.5772618977649913
$ws.Range("H4").Value = 0.7134762835499444
$ws.Range("I4").Value = 0.5147699734557172
$ws.Range("J4").Value = 0.112325968116517
$ws.Range("M4").Value = 0.376317243631263
$ws.Range("O4").Value = 2.559310836673859
$ws.Range("B5").Value = 0.6904324745417227
$ws.Range("C5").Value = 0.1616839603443623
$ws.Range("D5").Value = 0.5947118716952531
$ws.Range("E5").Value = 0.2307795469218448
$ws.Range("G5").Value = 0.5785243288502073
$ws.Range("H5").Value = 0.7153635194277399
$ws.Range("I5").Value = 0.5173516384041221
$ws.Range("J5").Value = 0.1126128126111556
$ws.Range("M5").Value = 0.3690585164700835
$ws.Range("O5").Value = 2.565852117372614
$ws.Range("B6").Value = 0.6861593960490779
$ws.Range("C6").Value = 0.1605698050372837
$ws.Range("D6").Value = 0.5946222955091969
$ws.Range("E6").Value = 0.230806649418966
$ws.Range("G6").Value = 0.5787416721767968
$ws.Range("H6").Value = 0.7156829311422825
$ws.Range("I6").Value = 0.5177862682578702
$ws.Range("J6").Value = 0.1126613194100283
$ws.Range("M6").Value = 0.3678542805294782
$ws.Range("O6").Value = 2.566967125426913
$ws.Range("B7").Value = 0.7158063008159559
$ws.Range("C7").Value = 0.1682980585337361
$ws.Range("D7").Value = 0.5952711396186601
$ws.Range("E7").Value = 0.2306291942717138
$ws.Range("G7").Value = 0.5772784057956528
$ws.Range("H7").Value = 0.7135013307006517
$ws.Range("I7").Value = 0.5148043919515519
$ws.Range("J7").Value = 0.1123297778399976
$ws.Range("M7").Value = 0.3762192784475431
$ws.Range("O7").Value = 2.559397121303263
$ws.Range("B8").Value = 0.8462999429476668
$ws.Range("C8").Value = 0.2022694871328383
$ws.Range("D8").Value = 0.5988040882707679
$ws.Range("E8").Value = 0.2301099816856187
$ws.Range("G8").Value = 0.5719450933739481
$ws.Range("H8").Value = 0.7047522295302073
$ws.Range("I8").Value = 0.5025115218931244
$ws.Range("J8").Value = 0.1109941968415953
$ws.Range("M8").Value = 0.4132846066916827
$ws.Range("O8").Value = 2.530188119146459
$ws.Range("B9").Value = 1.10090658011967
$ws.Range("C9").Value = 0.2683893455352688
$ws.Range("D9").Value = 0.6081352656146919
$ws.Range("E9").Value = 0.2300397551937472
$ws.Range("G9").Value = 0.5655448788371871
$ws.Range("H9").Value = 0.6907525623900455
$ws.Range("I9").Value = 0.4815268131907722
$ws.Range("J9").Value = 0.1088335348039848
$ws.Range("M9").Value = 0.4864917943567875
$ws.Range("O9").Value = 2.488035991959293
$ws.Range("B10").Value = 1.287114655970242
$ws.Range("C10").Value = 0.3166480035827703
$ws.Range("D10").Value = 0.6164543235250051
$ws.Range("E10").Value = 0.2305661824406258
$ws.Range("G10").Value = 0.5633357274210908
$ws.Range("H10").Value = 0.6823941256969874
$ws.Range("I10").Value = 0.4680217625935192
$ws.Range("J10").Value = 0.1075259273471012
$ws.Range("M10").Value = 0.5405798417600636
$ws.Range("O10").Value = 2.466334048115726
$ws.Range("B11").Value = 1.37162834635717
$ws.Range("C11").Value = 0.3385296678310965
$ws.Range("D11").Value = 0.6205564380545923
$ws.Range("E11").Value = 0.2309312356636184
$ws.Range("G11").Value = 0.5628761873916801
$ws.Range("H11").Value = 0.6790106438327825
$ws.Range("I11").Value = 0.4622952316467597
$ws.Range("J11").Value = 0.1069919101082668
$ws.Range("M11").Value = 0.5652487102281327
$ws.Range("O11").Value = 2.458483488048273
$ws.Range("B12").Value = 1.40360226415487
$ws.Range("C12").Value = 0.3468050553096305
$ws.Range("D12").Value = 0.6221554408774352
$ws.Range("E12").Value = 0.231087527777035
$ws.Range("G12").Value = 0.5627809197730329
$ws.Range("H12").Value = 0.6777896770854852
$ws.Range("I12").Value = 0.4601868834611587
$ws.Range("J12").Value = 0.1067984449505026
$ws.Range("M12").Value = 0.5745989981202371
$ws.Range("O12").Value = 2.455802188948354
$ws.Range("B13").Value = 1.396717452848407
$ws.Range("C13").Value = 0.3450232866324541
$ws.Range("D13").Value = 0.6218090390920565
$ws.Range("E13").Value = 0.2310530645594078
$ws.Range("G13").Value = 0.5627979298384389
$ws.Range("H13").Value = 0.6780499520066883
$ws.Range("I13").Value = 0.4606382748919948
$ws.Range("J13").Value = 0.106839721571319
$ws.Range("M13").Value = 0.5725848658633765
$ws.Range("O13").Value = 2.456366675601316
$ws.Range("B14").Value = 1.374259462876353
$ws.Range("C14").Value = 0.339210706431345
$ws.Range("D14").Value = 0.6206870750575604
$ws.Range("E14").Value = 0.2309437320942997
$ws.Range("G14").Value = 0.5628667695665968
$ws.Range("H14").Value = 0.6789089856331429
$ws.Range("I14").Value = 0.4621205696190707
$ws.Range("J14").Value = 0.1069758180828444
$ws.Range("M14").Value = 0.5660177919170479
$ws.Range("O14").Value = 2.458257048479652
$ws.Range("B15").Value = 1.36049939026725
$ws.Range("C15").Value = 0.3356489233451612
$ws.Range("D15").Value = 0.6200057788673803
$ws.Range("E15").Value = 0.2308791139061128
$ws.Range("G15").Value = 0.5629192009698869
$ws.Range("H15").Value = 0.6794430211108136
$ws.Range("I15").Value = 0.4630363591908591
$ws.Range("J15").Value = 0.1070603216508488
$ws.Range("M15").Value = 0.5619963915047066
$ws.Range("O15").Value = 2.459452947248678
$ws.Range("B16").Value = 1.281587444833065
$ws.Range("C16").Value = 0.3152165091655945
$ws.Range("D16").Value = 0.6161926284362664
$ws.Range("E16").Value = 0.2305448520042219
$ws.Range("G16").Value = 0.5633767619890193
$ws.Range("H16").Value = 0.682623676877526
$ws.Range("I16").Value = 0.4684044130066418
$ws.Range("J16").Value = 0.1075620513147442
$ws.Range("M16").Value = 0.5389689194650771
$ws.Range("O16").Value = 2.46688785309135
$ws.Range("B17").Value = 1.233126738108808
$ws.Range("C17").Value = 0.3026632651291266
$ws.Range("D17").Value = 0.6139347144653868
$ws.Range("E17").Value = 0.2303719542661327
$ws.Range("G17").Value = 0.5637973829269214
$ws.Range("H17").Value = 0.6846822166383362
$ws.Range("I17").Value = 0.4718044843600353
$ws.Range("J17").Value = 0.1078854305797563
$ws.Range("M17").Value = 0.5248583518348084
$ws.Range("O17").Value = 2.471967326915632
$ws.Range("B18").Value = 1.205235333107566
$ws.Range("C18").Value = 0.2954362734810161
$ws.Range("D18").Value = 0.6126659327363768
$ws.Range("E18").Value = 0.2302843283548128
$ws.Range("G18").Value = 0.5640906308668008
$ws.Range("H18").Value = 0.6859056513496995
$ws.Range("I18").Value = 0.4737993518949466
$ws.Range("J18").Value = 0.1080771538559802
$ws.Range("M18").Value = 0.5167483733287099
$ws.Range("O18").Value = 2.475079155577333
$ws.Range("B19").Value = 1.195788730634661
$ws.Range("C19").Value = 0.2929882017460557
$ws.Range("D19").Value = 0.6122414843200659
$ws.Range("E19").Value = 0.2302566899869412
$ws.Range("G19").Value = 0.5641987236935222
$ws.Range("H19").Value = 0.6863266534292336
$ws.Range("I19").Value = 0.4744815124706641
$ws.Range("J19").Value = 0.1081430508824894
$ws.Range("M19").Value = 0.5140035269406127
$ws.Range("O19").Value = 2.476165418108337
$ws.Range("B20").Value = 1.238287348123947
$ws.Range("C20").Value = 0.3040002758507114
$ws.Range("D20").Value = 0.6141719781181223
$ws.Range("E20").Value = 0.230389136196262
$ws.Range("G20").Value = 0.5637472935731722
$ws.Range("H20").Value = 0.6844590016125807
$ws.Range("I20").Value = 0.4714384786652133
$ws.Range("J20").Value = 0.1078504138141057
$ws.Range("M20").Value = 0.5263598234604245
$ws.Range("O20").Value = 2.471406912240781
$ws.Range("B21").Value = 1.380856735946224
$ws.Range("C21").Value = 0.3409182967512834
$ws.Range("D21").Value = 0.6210153854944451
$ws.Range("E21").Value = 0.230975355728237
$ws.Range("G21").Value = 0.5628444098666563
$ws.Range("H21").Value = 0.6786550300574135
$ws.Range("I21").Value = 0.4616835492550457
$ws.Range("J21").Value = 0.1069356055362256
$ws.Range("M21").Value = 0.5679464677379826
$ws.Range("O21").Value = 2.457693881452826
$ws.Range("B22").Value = 1.473860824318535
$ws.Range("C22").Value = 0.3649836441217644
$ws.Range("D22").Value = 0.6257538305883088
$ws.Range("E22").Value = 0.2314637171221321
$ws.Range("G22").Value = 0.5627134872083559
$ws.Range("H22").Value = 0.6752132104121387
$ws.Range("I22").Value = 0.4556589194348568
$ws.Range("J22").Value = 0.1063887642290418
$ws.Range("M22").Value = 0.595176369940404
$ws.Range("O22").Value = 2.450431301072825
$ws.Range("B23").Value = 1.424239230639614
$ws.Range("C23").Value = 0.3521454015338463
$ws.Range("D23").Value = 0.6232005276183656
$ws.Range("E23").Value = 0.2311934418435939
$ws.Range("G23").Value = 0.5627412441352533
$ws.Range("H23").Value = 0.6770180004903921
$ws.Range("I23").Value = 0.4588422143441804
$ws.Range("J23").Value = 0.1066759504405859
$ws.Range("M23").Value = 0.5806387816190863
$ws.Range("O23").Value = 2.454151689380723
$ws.Range("B24").Value = 1.235954331826974
$ws.Range("C24").Value = 0.3033958443096765
$ws.Range("D24").Value = 0.6140646198782065
$ws.Range("E24").Value = 0.2303813315666403
$ws.Range("G24").Value = 0.5637697788006477
$ws.Range("H24").Value = 0.684559792677419
$ws.Range("I24").Value = 0.471603824878688
$ws.Range("J24").Value = 0.1078662268019812
$ws.Range("M24").Value = 0.5256810007068538
$ws.Range("O24").Value = 2.471659679037771
$ws.Range("B25").Value = 1.032173665913831
$ws.Range("C25").Value = 0.2505571565302205
$ws.Range("D25").Value = 0.6053538285881501
$ws.Range("E25").Value = 0.2299572558731349
$ws.Range("G25").Value = 0.5668398599104449
$ws.Range("H25").Value = 0.6942015906986398
$ws.Range("I25").Value = 0.4868685724733215
$ws.Range("J25").Value = 0.1093689370378428
$ws.Range("M25").Value = 0.4666330331801163
$ws.Range("O25").Value = 2.497815121271884
